$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "C2" "320018701475"
Set-TextValue "C3" "320018701497"
Set-TextValue "C4" "320018701523"
Set-TextValue "C5" "320018701545"
Set-TextValue "D5" "320018701545"
Set-TextValue "C6" "320018701589"
Set-TextValue "D6" "320018701589"
Set-TextValue "C7" "320018701604"
Set-TextValue "D7" "320018701604"
Set-TextValue "C8" "320018701648"
Set-TextValue "C9" "320018702081"
Set-TextValue "C10" "320018702130"
Set-TextValue "C11" "320018702162"
Set-TextValue "C12" "320018702210"
Set-TextValue "C13" "320018702232"
Set-TextValue "D13" "320018702232"
Set-TextValue "C14" "320018702265"
Set-TextValue "D14" "320018702265"
Set-TextValue "C15" "320018702287"
Set-TextValue "D15" "320018702287"
Set-TextValue "C16" "320018702324"
Set-TextValue "D16" "320018702324"
Set-TextValue "C17" "320018702368"
Set-TextValue "D17" "320018702368"
Set-TextValue "C18" "320018702405"
Set-TextValue "C19" "320018702427"
Set-TextValue "C20" "320018702450"
Set-TextValue "C21" "320018702471"
Set-TextValue "C22" "320018702508"
